$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values (Terrence1 -> Terrence2)
$ws.Range("A2").Value = "Terrence2"
$ws.Range("B2").Value = "terrencereinhardt2@gmail.com"

# Remove rows 3 through 5 entirely (Kevin1, Laurencio1, Bro), shifting rows up
$ws.Rows("3:5").Delete()
